$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values are stored as text, matching the source
# data format (values like "1.006" or "211.10" must not be auto-converted
# to numbers, which would corrupt trailing zeros / multi-dot strings).
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '26.156.19'
$ws.Range("E2").Value = '  -1.08%  '

# Row 3
$ws.Range("D3").Value = '1.679.51'
$ws.Range("E3").Value = '  -0.74%  '

# Row 4
$ws.Range("D4").Value = '1.006'
$ws.Range("E4").Value = '  -0.44%  '

# Row 5
$ws.Range("D5").Value = '211.10'
$ws.Range("E5").Value = '  -3.62%  '

# Row 6
$ws.Range("D6").Value = '0.5293'
$ws.Range("E6").Value = '  -4.31%  '

# Row 7
$ws.Range("D7").Value = '1.006'
$ws.Range("E7").Value = '  -0.42%  '

# Row 8
$ws.Range("D8").Value = '0.2681'
$ws.Range("E8").Value = '  -1.53%  '

# Row 9
$ws.Range("D9").Value = '0.06304'
$ws.Range("E9").Value = '  -2.93%  '

# Row 10
$ws.Range("D10").Value = '21.29'
$ws.Range("E10").Value = '  -3.83%  '

# Row 11
$ws.Range("D11").Value = '0.07545'
$ws.Range("E11").Value = '  -0.75%  '

# Row 12
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.685.13'
$ws.Range("E12").Value = '  -0.47%  '

# Row 13
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = '4.498'
$ws.Range("E13").Value = '  -1.54%  '

# Row 14
$ws.Range("D14").Value = '0.5660'
$ws.Range("E14").Value = '  -3.32%  '

# Row 15
$ws.Range("D15").Value = '0.000008123'
$ws.Range("E15").Value = '  -4.13%  '

# Row 16
$ws.Range("D16").Value = '66.18'
$ws.Range("E16").Value = '  +1.08%  '

# Row 17
$ws.Range("D17").Value = '26.219.62'
$ws.Range("E17").Value = '  -1.15%  '

# Row 18
$ws.Range("B18").Value = 'Dai'
$ws.Range("C18").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D18").Value = '1.005'
$ws.Range("E18").Value = '  -0.44%  '

# Row 19
$ws.Range("B19").Value = 'Uniswap'
$ws.Range("C19").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D19").Value = '4.858'
$ws.Range("E19").Value = '  -2.33%  '

# Row 20
$ws.Range("D20").Value = '10.55'
$ws.Range("E20").Value = '  -4.04%  '

# Row 21
$ws.Range("D21").Value = '189.32'
$ws.Range("E21").Value = '  -0.84%  '

# Row 22
$ws.Range("D22").Value = '6.202'
$ws.Range("E22").Value = '  -1.05%  '

# Row 23
$ws.Range("D23").Value = '1.007'
$ws.Range("E23").Value = '  -0.39%  '

# Row 24
$ws.Range("D24").Value = '148.55'
$ws.Range("E24").Value = '  -1.00%  '

# Row 25
$ws.Range("D25").Value = '0.1261'
$ws.Range("E25").Value = '  -4.05%  '

# Row 26
$ws.Range("D26").Value = '7.625'
$ws.Range("E26").Value = '  -4.02%  '

# Row 27
$ws.Range("D27").Value = '15.86'
$ws.Range("E27").Value = '  +0.26%  '

# Row 28
$ws.Range("D28").Value = '0.06451'
$ws.Range("E28").Value = '  +1.94%  '

# Row 29
$ws.Range("D29").Value = '1.339'
$ws.Range("E29").Value = '  -5.09%  '

# Row 30
$ws.Range("D30").Value = '1.286'
$ws.Range("E30").Value = '  -3.41%  '

# Row 31
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").Value = '3.542'
$ws.Range("E31").Value = '  -1.43%  '

# Row 32
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = '3.477'
$ws.Range("E32").Value = '  -3.21%  '

# Row 33
$ws.Range("D33").Value = '1.653'
$ws.Range("E33").Value = '  -1.66%  '

# Row 34
$ws.Range("D34").Value = '1.008'
$ws.Range("E34").Value = '  -3.84%  '

# Row 35
$ws.Range("D35").Value = '0.6103'
$ws.Range("E35").Value = '  -2.55%  '

# Row 36
$ws.Range("D36").Value = '2.418'
$ws.Range("E36").Value = '  +0.62%  '

# Row 37
$ws.Range("D37").Value = '2.714'
$ws.Range("E37").Value = '  -0.25%  '

# Row 38
$ws.Range("D38").Value = '6.183'
$ws.Range("E38").Value = '  -1.10%  '

# Row 39
$ws.Range("B39").Value = 'Maker'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D39").Value = '1.098.92'
$ws.Range("E39").Value = '  -2.39%  '

# Row 40
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").Value = '0.01608'
$ws.Range("E40").Value = '  -2.38%  '

# Row 41
$ws.Range("D41").Value = '0.8703'
$ws.Range("E41").Value = '  -1.63%  '

# Row 42
$ws.Range("D42").Value = '1.006'
$ws.Range("E42").Value = '  -0.96%  '

# Row 43
$ws.Range("D43").Value = '100.00'
$ws.Range("E43").Value = '  -0.84%  '

# Row 44
$ws.Range("D44").Value = '1.834.54'
$ws.Range("E44").Value = '  -0.45%  '

# Row 45
$ws.Range("D45").Value = '0.00000000109'
$ws.Range("E45").Value = '  -5.33%  '

# Row 46
$ws.Range("D46").Value = '57.00'
$ws.Range("E46").Value = '  -1.16%  '

# Row 47
$ws.Range("D47").Value = '1.004'
$ws.Range("E47").Value = '  -0.46%  '

# Row 48
$ws.Range("D48").Value = '0.05288'
$ws.Range("E48").Value = '  +0.11%  '

# Row 49
$ws.Range("D49").Value = '7.952'
$ws.Range("E49").Value = '  -3.75%  '

# Row 50
$ws.Range("E50").Value = '  -0.62%  '

# Row 51
$ws.Range("D51").Value = '5.959'
$ws.Range("E51").Value = '  -2.40%  '
